$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Total Households values for rows 2-98 (column D), corresponding to the
# walkshed_proportional_households_stations_p75 output reconciliation.
$values = @(4450,6579,12975,3154,4375,19982,12549,13018,14300,5806,8011,7206,19587,2680,21890,13226,3577,37923,11315,24223,14726,7243,7231,37270,5942,23404,11402,31718,26641,17479,9216,17383,5026,11253,3511,11594,21133,21454,5885,4809,8579,9329,11409,4823,6310,15857,12001,12296,6616,5214,0,7475,22876,6706,20512,8962,5487,35583,20929,8424,3844,24616,11680,16518,13024,8342,12288,7046,9846,16688,9193,5664,39422,18447,4823,8580,7703,11445,6177,8615,10413,9165,10119,43023,24154,11056,13197,7653,21478,18794,4393,8744,7637,9914,6754,19200,0)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}
